$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3450.5557
$ws.Range("I132").Value = 2794.6428
$ws.Range("J132").Value = 5746.25
$ws.Range("K132").Value = 8383.928400000001
$ws.Range("L132").Value = 17238.75
$ws.Range("M132").Value = -5853.928400000001
$ws.Range("N132").Value = -22298.75
$ws.Range("H133").Value = 30118.75
$ws.Range("J133").Value = 30118.75
$ws.Range("L133").Value = 30118.75
$ws.Range("N133").Value = -40238.75
$ws.Range("H134").Value = 30923.076
$ws.Range("J134").Value = 30923.076
$ws.Range("L134").Value = 30923.076
$ws.Range("N134").Value = -41063.076
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null
$ws.Range("H137").Value = 1706.6786
$ws.Range("I137").Value = 1495.6
$ws.Range("J137").Value = 1752.5652
$ws.Range("K137").Value = 4486.799999999999
$ws.Range("L137").Value = 5257.6956
$ws.Range("M137").Value = -1936.799999999999
$ws.Range("N137").Value = -10357.6956
$ws.Range("H138").Value = 3485.1292
$ws.Range("I138").Value = 1433.2916
$ws.Range("J138").Value = 4781.0264
$ws.Range("K138").Value = 4299.8748
$ws.Range("L138").Value = 14343.0792
$ws.Range("M138").Value = 840.1252000000004
$ws.Range("N138").Value = -24623.0792
$ws.Range("H140").Value = 34666.668
$ws.Range("J140").Value = 34666.668
$ws.Range("L140").Value = 34666.668
$ws.Range("N140").Value = -45026.668
$ws.Range("H141").Value = 5919
$ws.Range("I141").Value = 5648.75
$ws.Range("K141").Value = 16946.25
$ws.Range("M141").Value = -11766.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3735.2727
$ws.Range("I32").Value = 3035.7097
$ws.Range("J32").Value = 6626.8
$ws.Range("K32").Value = 3035.7097
$ws.Range("L32").Value = 6626.8
$ws.Range("M32").Value = -2748.7097
$ws.Range("N32").Value = -7200.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 742.7778
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 812.1429000000001
$ws.Range("K80").Value = 500
$ws.Range("L80").Value = 812.1429000000001
$ws.Range("M80").Value = 498
$ws.Range("N80").Value = -2808.1429
$ws.Range("H83").Value = 742.7778
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 812.1429000000001
$ws.Range("K83").Value = 2500
$ws.Range("L83").Value = 4060.7145
$ws.Range("M83").Value = 2492
$ws.Range("N83").Value = -14044.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 10000137
$ws.Range("J23").Value = 244.8
$ws.Range("L23").Value = 734.4000000000001
$ws.Range("N23").Value = -1204.4
$ws.Range("H70").Value = 151202.88
$ws.Range("I70").Value = 151202.88
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 453608.64
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -453293.64
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 151202.88
$ws.Range("I73").Value = 151202.88
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 453608.64
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -452516.64
$ws.Range("N73").Value = $null
$ws.Range("H88").Value = 3788.889
$ws.Range("J88").Value = 3788.889
$ws.Range("L88").Value = 11366.667
$ws.Range("N88").Value = -12222.667
$ws.Range("H91").Value = 3788.889
$ws.Range("J91").Value = 3788.889
$ws.Range("L91").Value = 11366.667
$ws.Range("N91").Value = -14330.667
$ws.Range("H139").Value = 86776.08
$ws.Range("I139").Value = 123432.11
$ws.Range("J139").Value = 4300
$ws.Range("K139").Value = 370296.33
$ws.Range("L139").Value = 12900
$ws.Range("M139").Value = -365156.33
$ws.Range("N139").Value = -23180
$ws.Range("H140").Value = 15498.1
$ws.Range("I140").Value = 15498.1
$ws.Range("K140").Value = 46494.3
$ws.Range("M140").Value = -41314.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5963.114
$ws.Range("I70").Value = 5993.4287
$ws.Range("J70").Value = 5845.222
$ws.Range("K70").Value = 5993.4287
$ws.Range("L70").Value = 5845.222
$ws.Range("M70").Value = -5723.4287
$ws.Range("N70").Value = -6385.222
$ws.Range("H73").Value = 5963.114
$ws.Range("I73").Value = 5993.4287
$ws.Range("J73").Value = 5845.222
$ws.Range("K73").Value = 5993.4287
$ws.Range("L73").Value = 5845.222
$ws.Range("M73").Value = -5057.4287
$ws.Range("N73").Value = -7717.222
$ws.Range("H80").Value = 2558.8462
$ws.Range("I80").Value = 2547.8
$ws.Range("J80").Value = 2578.5715
$ws.Range("K80").Value = 2547.8
$ws.Range("L80").Value = 2578.5715
$ws.Range("M80").Value = -1549.8
$ws.Range("N80").Value = -4574.5715
$ws.Range("H83").Value = 2558.8462
$ws.Range("I83").Value = 2547.8
$ws.Range("J83").Value = 2578.5715
$ws.Range("K83").Value = 12739
$ws.Range("L83").Value = 12892.8575
$ws.Range("M83").Value = -7747
$ws.Range("N83").Value = -22876.8575
$ws.Range("H132").Value = 4220.9116
$ws.Range("I132").Value = 7700
$ws.Range("J132").Value = 3621.0688
$ws.Range("K132").Value = 23100
$ws.Range("L132").Value = 10863.2064
$ws.Range("M132").Value = -20570
$ws.Range("N132").Value = -15923.2064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 916
$ws.Range("I16").Value = 916
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 916
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -746
$ws.Range("N16").Value = $null
$ws.Range("H46").Value = 33334370
$ws.Range("I46").Value = 37038076
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 37038076
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -37037888
$ws.Range("N46").Value = -1376
$ws.Range("H139").Value = 63290.668
$ws.Range("J139").Value = 63290.668
$ws.Range("L139").Value = 63290.668
$ws.Range("N139").Value = -73570.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2500.8
$ws.Range("I96").Value = 2666.6667
$ws.Range("J96").Value = 2252
$ws.Range("K96").Value = 2666.6667
$ws.Range("L96").Value = 2252
$ws.Range("M96").Value = -1293.6667
$ws.Range("N96").Value = -4998
